# fixede accuracy til test data istedet for train data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3's Scheduler ("None") was actually the data meant for the new final
# row; clear it here (it becomes an empty inline string) and re-add "None"
# on the newly appended row below.
$ws.Range("T3").Value = ""

# New run logged at 2024-1-6 14:16:55 (row 4)
$ws.Range("A4").Value = "2024-1-6 14:16:55"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 64
$ws.Range("D4").Value = 0.001
$ws.Range("E4").Value = "ADAM"
$ws.Range("F4").Value = "CEL"
$ws.Range("G4").Value = 36.4
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 0.6728
$ws.Range("J4").Value = 0.5442
$ws.Range("K4").Value = 0.5196433546949011
$ws.Range("L4").Value = "FER2013"
$ws.Range("M4").Value = "cpu"
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = "Alfred"
$ws.Range("Q4").Value = 364.4
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = 0

# New run logged at 2024-1-6 14:32:46 (row 5)
$ws.Range("A5").Value = "2024-1-6 14:32:46"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 64
$ws.Range("D5").Value = 0.001
$ws.Range("E5").Value = "ADAM"
$ws.Range("F5").Value = "CEL"
$ws.Range("G5").Value = 36.9
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 1.0813
$ws.Range("J5").Value = 1.0813
$ws.Range("K5").Value = 0.515881861242686
$ws.Range("L5").Value = "FER2013"
$ws.Range("M5").Value = "cpu"
$ws.Range("N5").Value = 4
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = "Alfred"
$ws.Range("Q5").Value = 369.3
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.005
$ws.Range("T5").Value = "None"
$ws.Range("U5").Value = 0
